$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lowercase the "Desde"/"Hasta" header labels in D1/E1
$ws.Range("D1").Value = "desde"
$ws.Range("E1").Value = "hasta"

# Update the saved selection/active cell from E3 to E4
$ws.Range("E4").Select()
